$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.592.69'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.701.87'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.93'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3922'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4069'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.490'
$ws.Range("E9").Value = '  -2.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9983'
$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.45'
$ws.Range("E11").Value = '  +0.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08832'
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.30'
$ws.Range("E13").Value = '  +11.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.480'
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.127'
$ws.Range("E15").Value = '  +0.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001361'
$ws.Range("E16").Value = '  +3.30%  '

$ws.Range("D17").Value = '1.702.94'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.77'
$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07206'
$ws.Range("E19").Value = '  +2.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.54'
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.311'
$ws.Range("E21").Value = '  +3.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").Value = '  -2.25%  '

$ws.Range("D24").Value = '24.594.25'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.025'
$ws.Range("E25").Value = '  -3.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.333'
$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.08'
$ws.Range("E27").Value = '  +2.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.23'
$ws.Range("E28").Value = '  +2.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.973'
$ws.Range("E29").Value = '  +16.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '144.32'
$ws.Range("E30").Value = '  +6.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.397'
$ws.Range("E31").Value = '  -4.49%  '

$ws.Range("D32").Value = '1.884.68'
$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08816'
$ws.Range("E33").Value = '  -2.18%  '

$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.167'
$ws.Range("E34").Value = '  +10.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.059'
$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.207'
$ws.Range("E36").Value = '  -5.31%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.03126'
$ws.Range("E37").Value = '  +6.65%  '

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8794'
$ws.Range("E38").Value = '  +14.92%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2812'
$ws.Range("E39").Value = '  +2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.91'
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09199'
$ws.Range("E41").Value = '  +0.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.27'
$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.483'
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.31'
$ws.Range("E44").Value = '  +7.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7541'
$ws.Range("E45").Value = '  +5.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.678'
$ws.Range("E46").Value = '  +3.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.260'
$ws.Range("E47").Value = '  +1.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.403'
$ws.Range("E48").Value = '  +5.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9988'
$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.90'
$ws.Range("E50").Value = '  +0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08255'
$ws.Range("E51").Value = '  +3.54%  '
